# Updated cryptos list values (Price / Volume(1h)) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.Value = "'62.281.66"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "

$dCell = $ws.Range("D3")
$dCell.Value = "'2.448.99"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.11%  "

$dCell = $ws.Range("D5")
$dCell.Value = "'583.27"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +3.58%  "

$dCell = $ws.Range("D6")
$dCell.Value = "'143.22"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$dCell = $ws.Range("D8")
$dCell.Value = "'0.532"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +0.84%  "

$dCell = $ws.Range("D9")
$dCell.Value = "'2.443.61"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("E11").Value = "  +2.68%  "

$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("E13").Value = "  -2.07%  "

$dCell = $ws.Range("D14")
$dCell.Value = "'26.48"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("E15").Value = "  +3.00%  "

$ws.Range("E16").Value = "  +0.86%  "

$dCell = $ws.Range("D17")
$dCell.Value = "'62.165.99"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "

$dCell = $ws.Range("D18")
$dCell.Value = "'2.440.19"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "

$dCell = $ws.Range("D19")
$dCell.Value = "'10.81"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  -2.02%  "

$dCell = $ws.Range("D20")
$dCell.Value = "'7.11"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "

$dCell = $ws.Range("D21")
$dCell.Value = "'327.64"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -3.27%  "

$dCell = $ws.Range("D25")
$dCell.Value = "'65.77"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "

$dCell = $ws.Range("D26")
$dCell.Value = "'9.22"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +3.59%  "

$dCell = $ws.Range("D27")
$dCell.Value = "'593.53"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -4.52%  "

$dCell = $ws.Range("D28")
$dCell.Value = "'0.0₃0975"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +2.30%  "

$dCell = $ws.Range("D29")
$dCell.Value = "'2.570.55"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +0.82%  "

$dCell = $ws.Range("D30")
$dCell.Value = "'0.994"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -0.82%  "

$ws.Range("E31").Value = "  -1.30%  "

$dCell = $ws.Range("D32")
$dCell.Value = "'7.99"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("E33").Value = "  +2.34%  "

$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("E35").Value = "  -1.93%  "

$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("E37").Value = "  -1.35%  "

$dCell = $ws.Range("D38")
$dCell.Value = "'0.378"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "

$dCell = $ws.Range("D39")
$dCell.Value = "'153.20"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +4.53%  "

$dCell = $ws.Range("D40")
$dCell.Value = "'18.43"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -1.22%  "

$dCell = $ws.Range("D41")
$dCell.Value = "'5.28"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "

$dCell = $ws.Range("D42")
$dCell.Value = "'42.86"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +0.63%  "

$dCell = $ws.Range("D43")
$dCell.Value = "'1.71"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -1.35%  "

$ws.Range("E44").Value = "  +0.03%  "

$dCell = $ws.Range("D45")
$dCell.Value = "'2.52"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +2.28%  "

$dCell = $ws.Range("D46")
$dCell.Value = "'142.55"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -1.84%  "

$ws.Range("E47").Value = "  -0.82%  "

$dCell = $ws.Range("D48")
$dCell.Value = "'0.0₆0258"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +18.16%  "

$ws.Range("E49").Value = "  +2.23%  "

$dCell = $ws.Range("D50")
$dCell.Value = "'0.0523"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "

$dCell = $ws.Range("D51")
$dCell.Value = "'19.91"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "

